$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 45.98144433333334
$ws.Range("H2").Value = 137.944333
$ws.Range("I2").Value = 0.9841234286873372
$ws.Range("J2").Value = 0.984123428687337
$ws.Range("M2").Value = 61.156892
$ws.Range("N2").Value = 183.470676
$ws.Range("O2").Value = 0.9308124486389074
$ws.Range("P2").Value = 0.9308124486389074
$ws.Range("Q2").Value = 2812.082225097679
$ws.Range("R2").Value = 25308.74002587911
$ws.Range("S2").Value = 0.9160343384193775
$ws.Range("T2").Value = 0.9160343384193774

$ws.Range("G3").Value = 45.98144433333334
$ws.Range("H3").Value = 137.944333
$ws.Range("I3").Value = 0.9841234286873372
$ws.Range("J3").Value = 0.984123428687337
$ws.Range("O3").Value = 0.02171808228502914
$ws.Range("P3").Value = 0.02171808228502914
$ws.Range("Q3").Value = 65.61260890552555
$ws.Range("R3").Value = 590.51348014973
$ws.Range("S3").Value = 0.0213732736028566
$ws.Range("T3").Value = 0.02137327360285659

$ws.Range("G4").Value = 45.98144433333334
$ws.Range("H4").Value = 137.944333
$ws.Range("I4").Value = 0.9841234286873372
$ws.Range("J4").Value = 0.984123428687337
$ws.Range("M4").Value = 2.00294
$ws.Range("N4").Value = 6.00882
$ws.Range("O4").Value = 0.03048489589491914
$ws.Range("P4").Value = 0.03048489589491914
$ws.Range("Q4").Value = 92.09807411300667
$ws.Range("R4").Value = 828.88266701706
$ws.Range("S4").Value = 0.03000090027128436
$ws.Range("T4").Value = 0.03000090027128435

$ws.Range("G5").Value = 45.98144433333334
$ws.Range("H5").Value = 137.944333
$ws.Range("I5").Value = 0.9841234286873372
$ws.Range("J5").Value = 0.984123428687337
$ws.Range("M5").Value = 1.115932333333334
$ws.Range("N5").Value = 3.347797
$ws.Range("O5").Value = 0.01698457318114416
$ws.Range("P5").Value = 0.01698457318114415
$ws.Range("Q5").Value = 51.31218046493345
$ws.Range("R5").Value = 461.8096241844011
$ws.Range("S5").Value = 0.01671491639381858
$ws.Range("T5").Value = 0.01671491639381858

$ws.Range("I6").Value = 0.002244435796517234
$ws.Range("J6").Value = 0.002244435796517234
$ws.Range("M6").Value = 61.156892
$ws.Range("N6").Value = 183.470676
$ws.Range("O6").Value = 0.9308124486389074
$ws.Range("P6").Value = 0.9308124486389074
$ws.Range("Q6").Value = 6.413360178994666
$ws.Range("R6").Value = 57.72024161095199
$ws.Range("S6").Value = 0.002089148779569024
$ws.Range("T6").Value = 0.002089148779569023

$ws.Range("I7").Value = 0.002244435796517234
$ws.Range("J7").Value = 0.002244435796517234
$ws.Range("O7").Value = 0.02171808228502914
$ws.Range("P7").Value = 0.02171808228502914
$ws.Range("S7").Value = 0.00004874484131222622
$ws.Range("T7").Value = 0.0000487448413122262

$ws.Range("I8").Value = 0.002244435796517234
$ws.Range("J8").Value = 0.002244435796517234
$ws.Range("M8").Value = 2.00294
$ws.Range("N8").Value = 6.00882
$ws.Range("O8").Value = 0.03048489589491914
$ws.Range("P8").Value = 0.03048489589491914
$ws.Range("Q8").Value = 0.2100429766266667
$ws.Range("R8").Value = 1.89038678964
$ws.Range("S8").Value = 0.00006842139159965783
$ws.Range("T8").Value = 0.0000684213915996578

$ws.Range("I9").Value = 0.002244435796517234
$ws.Range("J9").Value = 0.002244435796517234
$ws.Range("M9").Value = 1.115932333333334
$ws.Range("N9").Value = 3.347797
$ws.Range("O9").Value = 0.01698457318114416
$ws.Range("P9").Value = 0.01698457318114415
$ws.Range("Q9").Value = 0.1170248479771111
$ws.Range("R9").Value = 1.053223631794
$ws.Range("S9").Value = 0.00003812078403632655
$ws.Range("T9").Value = 0.00003812078403632653

$ws.Range("G10").Value = 0.547937
$ws.Range("H10").Value = 1.643811
$ws.Range("I10").Value = 0.01172728797372169
$ws.Range("J10").Value = 0.01172728797372169
$ws.Range("M10").Value = 61.156892
$ws.Range("N10").Value = 183.470676
$ws.Range("O10").Value = 0.9308124486389074
$ws.Range("P10").Value = 0.9308124486389074
$ws.Range("Q10").Value = 33.510123931804
$ws.Range("R10").Value = 301.591115386236
$ws.Range("S10").Value = 0.0109159056347135
$ws.Range("T10").Value = 0.0109159056347135

$ws.Range("G11").Value = 0.547937
$ws.Range("H11").Value = 1.643811
$ws.Range("I11").Value = 0.01172728797372169
$ws.Range("J11").Value = 0.01172728797372169
$ws.Range("O11").Value = 0.02171808228502914
$ws.Range("P11").Value = 0.02171808228502914
$ws.Range("Q11").Value = 0.7818713963233332
$ws.Range("R11").Value = 7.036842566909999
$ws.Range("S11").Value = 0.0002546942051935204
$ws.Range("T11").Value = 0.0002546942051935203

$ws.Range("G12").Value = 0.547937
$ws.Range("H12").Value = 1.643811
$ws.Range("I12").Value = 0.01172728797372169
$ws.Range("J12").Value = 0.01172728797372169
$ws.Range("M12").Value = 2.00294
$ws.Range("N12").Value = 6.00882
$ws.Range("O12").Value = 0.03048489589491914
$ws.Range("P12").Value = 0.03048489589491914
$ws.Range("Q12").Value = 1.09748493478
$ws.Range("R12").Value = 9.877364413019999
$ws.Range("S12").Value = 0.0003575051530086431
$ws.Range("T12").Value = 0.000357505153008643

$ws.Range("G13").Value = 0.547937
$ws.Range("H13").Value = 1.643811
$ws.Range("I13").Value = 0.01172728797372169
$ws.Range("J13").Value = 0.01172728797372169
$ws.Range("M13").Value = 1.115932333333334
$ws.Range("N13").Value = 3.347797
$ws.Range("O13").Value = 0.01698457318114416
$ws.Range("P13").Value = 0.01698457318114415
$ws.Range("Q13").Value = 0.6114606149296667
$ws.Range("R13").Value = 5.503145534367
$ws.Range("S13").Value = 0.0001991829808060279
$ws.Range("T13").Value = 0.0001991829808060278

$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.08900066666666666
$ws.Range("H14").Value = 0.267002
$ws.Range("I14").Value = 0.001904847542424061
$ws.Range("J14").Value = 0.001904847542424061
$ws.Range("M14").Value = 61.156892
$ws.Range("N14").Value = 183.470676
$ws.Range("O14").Value = 0.9308124486389074
$ws.Range("P14").Value = 0.9308124486389074
$ws.Range("Q14").Value = 5.443004159261333
$ws.Range("R14").Value = 48.98703743335199
$ws.Range("S14").Value = 0.001773055805247546
$ws.Range("T14").Value = 0.001773055805247546

$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.08900066666666666
$ws.Range("H15").Value = 0.267002
$ws.Range("I15").Value = 0.001904847542424061
$ws.Range("J15").Value = 0.001904847542424061
$ws.Range("O15").Value = 0.02171808228502914
$ws.Range("P15").Value = 0.02171808228502914
$ws.Range("Q15").Value = 0.1269983146244444
$ws.Range("R15").Value = 1.14298483162
$ws.Range("S15").Value = 0.00004136963566680131
$ws.Range("T15").Value = 0.0000413696356668013

$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.08900066666666666
$ws.Range("H16").Value = 0.267002
$ws.Range("I16").Value = 0.001904847542424061
$ws.Range("J16").Value = 0.001904847542424061
$ws.Range("M16").Value = 2.00294
$ws.Range("N16").Value = 6.00882
$ws.Range("O16").Value = 0.03048489589491914
$ws.Range("P16").Value = 0.03048489589491914
$ws.Range("Q16").Value = 0.1782629952933333
$ws.Range("R16").Value = 1.60436695764
$ws.Range("S16").Value = 0.00005806907902649009
$ws.Range("T16").Value = 0.00005806907902649008

$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.08900066666666666
$ws.Range("H17").Value = 0.267002
$ws.Range("I17").Value = 0.001904847542424061
$ws.Range("J17").Value = 0.001904847542424061
$ws.Range("M17").Value = 1.115932333333334
$ws.Range("N17").Value = 3.347797
$ws.Range("O17").Value = 0.01698457318114416
$ws.Range("P17").Value = 0.01698457318114415
$ws.Range("Q17").Value = 0.09931872162155557
$ws.Range("R17").Value = 0.893868494594
$ws.Range("S17").Value = 0.00004874484131222622
$ws.Range("T17").Value = 0.0000487448413122262
